# Automatic update of files.
#
# 1) The "Förändrad" column (C) is refreshed from 46066 to 46070 for every
#    data row (rows 2-38).
# 2) Three groups of rows get re-sorted (their full contents, A:Z, trade
#    places) because the underlying data source re-ordered them:
#       - row 16  <->  row 17
#       - row 20  <->  row 21
#       - rows 32,33,34,36,37 rotate: 36->32, 37->33, 32->34, 33->36, 34->37
#         (row 35 stays in place)

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------
# Step 1: update the "Förändrad" date-serial column (C) for all rows.
# ---------------------------------------------------------------------
$ws.Range("C2:C38").Value = 46070

# ---------------------------------------------------------------------
# Step 2: swap full row contents (A:Z) where the row order changed.
# Using .Formula captures both literal values (numbers/dates/text) and
# formulas (HYPERLINK(...)) uniformly, so a plain swap reproduces the
# exact target state, including cells that only exist on one side
# (e.g. column Z).
# ---------------------------------------------------------------------

function Swap-Rows($rowA, $rowB) {
    $rangeA = $ws.Range("A$rowA`:Z$rowA")
    $rangeB = $ws.Range("A$rowB`:Z$rowB")
    $valA = $rangeA.Formula
    $valB = $rangeB.Formula
    $rangeA.Formula = $valB
    $rangeB.Formula = $valA
}

# row 16 <-> row 17  (A 54167-2024  <->  A 60875-2025)
Swap-Rows 16 17

# row 20 <-> row 21  (A 2864-2026  <->  A 61963-2025)
Swap-Rows 20 21

# rows 32-37 rotate: capture all five affected rows first (row 35 is
# unaffected and stays untouched), then write them to their new homes.
$row32 = $ws.Range("A32:Z32").Formula
$row33 = $ws.Range("A33:Z33").Formula
$row34 = $ws.Range("A34:Z34").Formula
$row36 = $ws.Range("A36:Z36").Formula
$row37 = $ws.Range("A37:Z37").Formula

$ws.Range("A32:Z32").Formula = $row36   # A 62183-2025 -> row 32
$ws.Range("A33:Z33").Formula = $row37   # A 5733-2026  -> row 33
$ws.Range("A34:Z34").Formula = $row32   # A 63676-2023 -> row 34
$ws.Range("A36:Z36").Formula = $row33   # A 13040-2024 -> row 36
$ws.Range("A37:Z37").Formula = $row34   # A 13384-2023 -> row 37
